# Auto-generated from diff: update computed market-price columns (H-N)
# across multiple worksheets (scheduled market-data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 20.428572
$ws.Range("I11").Value = 20.428572
$ws.Range("K11").Value = 20.428572
$ws.Range("M11").Value = 119.571428
$ws.Range("H17").Value = 1922.8667
$ws.Range("J17").Value = 1922.8667
$ws.Range("L17").Value = 5768.6001
$ws.Range("N17").Value = -6104.6001
$ws.Range("H19").Value = 1241.5555
$ws.Range("I19").Value = 1497.5
$ws.Range("J19").Value = 1036.8
$ws.Range("K19").Value = 1497.5
$ws.Range("L19").Value = 1036.8
$ws.Range("M19").Value = -1322.5
$ws.Range("N19").Value = -1386.8
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H55").Value = 821.2
$ws.Range("I55").Value = 537
$ws.Range("J55").Value = 1247.5
$ws.Range("K55").Value = 537
$ws.Range("L55").Value = 1247.5
$ws.Range("M55").Value = -323
$ws.Range("N55").Value = -1675.5
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H99").Value = 380.33334
$ws.Range("I99").Value = 284.33334
$ws.Range("J99").Value = 476.33334
$ws.Range("K99").Value = 853.0000200000001
$ws.Range("L99").Value = 1429.00002
$ws.Range("M99").Value = 644.9999799999999
$ws.Range("N99").Value = -4425.000019999999
$ws.Range("H101").Value = 33333748
$ws.Range("I101").Value = 50000120
$ws.Range("K101").Value = 150000360
$ws.Range("M101").Value = -149998738
$ws.Range("H113").Value = 2382.2856
$ws.Range("I113").Value = 2017.4
$ws.Range("J113").Value = 3294.5
$ws.Range("K113").Value = 2017.4
$ws.Range("L113").Value = 3294.5
$ws.Range("M113").Value = 1236.6
$ws.Range("N113").Value = -9802.5
$ws.Range("H137").Value = 1125.7778
$ws.Range("I137").Value = 1125.7778
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 3377.3334
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -827.3334000000004
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 4220.5312
$ws.Range("J138").Value = 5277.75
$ws.Range("L138").Value = 15833.25
$ws.Range("N138").Value = -26113.25
$ws.Range("H141").Value = 5856.143
$ws.Range("I141").Value = 6331
$ws.Range("J141").Value = 5500
$ws.Range("K141").Value = 18993
$ws.Range("L141").Value = 16500
$ws.Range("M141").Value = -13813
$ws.Range("N141").Value = -26860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 1495.5555
$ws.Range("J22").Value = 1495
$ws.Range("L22").Value = 1495
$ws.Range("N22").Value = -2093
$ws.Range("H32").Value = 2131
$ws.Range("I32").Value = 2085.225
$ws.Range("K32").Value = 2085.225
$ws.Range("M32").Value = -1798.225
$ws.Range("H46").Value = 2943.8
$ws.Range("I46").Value = 3405.6667
$ws.Range("J46").Value = 2251
$ws.Range("K46").Value = 3405.6667
$ws.Range("L46").Value = 2251
$ws.Range("M46").Value = -3086.6667
$ws.Range("N46").Value = -2889
$ws.Range("H88").Value = 900
$ws.Range("J88").Value = 900
$ws.Range("L88").Value = 900
$ws.Range("N88").Value = -1712
$ws.Range("H91").Value = 900
$ws.Range("J91").Value = 900
$ws.Range("L91").Value = 900
$ws.Range("N91").Value = -3708

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 66324.664
$ws.Range("J52").Value = 66324.664
$ws.Range("L52").Value = 66324.664
$ws.Range("N52").Value = -66850.664
$ws.Range("H75").Value = 79999.5
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 79999.5
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 79999.5
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -81871.5
$ws.Range("H78").Value = 79999.5
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 79999.5
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 239998.5
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -249358.5
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H121").Value = 66324.664
$ws.Range("J121").Value = 66324.664
$ws.Range("L121").Value = 66324.664
$ws.Range("N121").Value = -69818.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1884.8334
$ws.Range("I31").Value = 1312
$ws.Range("J31").Value = 1999.4
$ws.Range("K31").Value = 1312
$ws.Range("L31").Value = 1999.4
$ws.Range("M31").Value = -1017
$ws.Range("N31").Value = -2589.4
$ws.Range("H34").Value = 1884.8334
$ws.Range("I34").Value = 1312
$ws.Range("J34").Value = 1999.4
$ws.Range("K34").Value = 1312
$ws.Range("L34").Value = 1999.4
$ws.Range("M34").Value = -1110
$ws.Range("N34").Value = -2403.4
$ws.Range("H58").Value = 979.9666999999999
$ws.Range("I58").Value = 919.96
$ws.Range("K58").Value = 919.96
$ws.Range("M58").Value = -716.96
$ws.Range("H62").Value = 3084
$ws.Range("I62").Value = 2953.3
$ws.Range("K62").Value = 2953.3
$ws.Range("M62").Value = -2329.3
$ws.Range("H65").Value = 3084
$ws.Range("I65").Value = 2953.3
$ws.Range("K65").Value = 14766.5
$ws.Range("M65").Value = -11646.5
$ws.Range("H86").Value = 9667.333000000001
$ws.Range("I86").Value = 8999.5
$ws.Range("K86").Value = 8999.5
$ws.Range("M86").Value = -7876.5
$ws.Range("H89").Value = 9667.333000000001
$ws.Range("I89").Value = 8999.5
$ws.Range("K89").Value = 44997.5
$ws.Range("M89").Value = -39381.5
$ws.Range("H92").Value = 49495.5
$ws.Range("J92").Value = 49495.5
$ws.Range("L92").Value = 49495.5
$ws.Range("N92").Value = -54487.5
$ws.Range("H99").Value = 3750
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H126").Value = 3750
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H136").Value = 979.9666999999999
$ws.Range("I136").Value = 919.96
$ws.Range("K136").Value = 2759.88
$ws.Range("M136").Value = -209.8800000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H107").Value = 1363.7693
$ws.Range("I107").Value = 672.8570999999999
$ws.Range("K107").Value = 672.8570999999999
$ws.Range("M107").Value = 1247.1429
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 3080.9092
$ws.Range("I122").Value = 3027.4285
$ws.Range("K122").Value = 9082.2855
$ws.Range("M122").Value = -6632.2855
$ws.Range("H126").Value = 2855.375
$ws.Range("I126").Value = 2855.375
$ws.Range("K126").Value = 8566.125
$ws.Range("M126").Value = -6096.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 699
$ws.Range("I16").Value = 699
$ws.Range("K16").Value = 699
$ws.Range("M16").Value = -529
$ws.Range("H22").Value = 1456
$ws.Range("I22").Value = 1329.8
$ws.Range("J22").Value = 1666.3334
$ws.Range("K22").Value = 1329.8
$ws.Range("L22").Value = 1666.3334
$ws.Range("M22").Value = -1034.8
$ws.Range("N22").Value = -2256.3334
$ws.Range("H27").Value = 1456
$ws.Range("I27").Value = 1329.8
$ws.Range("J27").Value = 1666.3334
$ws.Range("K27").Value = 1329.8
$ws.Range("L27").Value = 1666.3334
$ws.Range("M27").Value = -1222.8
$ws.Range("N27").Value = -1880.3334
$ws.Range("H46").Value = 3909.3125
$ws.Range("I46").Value = 1250
$ws.Range("J46").Value = 4795.75
$ws.Range("K46").Value = 1250
$ws.Range("L46").Value = 4795.75
$ws.Range("M46").Value = -1062
$ws.Range("N46").Value = -5171.75
$ws.Range("H68").Value = 2700
$ws.Range("I68").Value = 2728.5715
$ws.Range("K68").Value = 2728.5715
$ws.Range("M68").Value = -1979.5715
$ws.Range("H71").Value = 2700
$ws.Range("I71").Value = 2728.5715
$ws.Range("K71").Value = 13642.8575
$ws.Range("M71").Value = -9898.8575

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H62").Value = 12540.1
$ws.Range("J62").Value = 10056.857
$ws.Range("L62").Value = 10056.857
$ws.Range("N62").Value = -11304.857
$ws.Range("H65").Value = 12540.1
$ws.Range("J65").Value = 10056.857
$ws.Range("L65").Value = 50284.285
$ws.Range("N65").Value = -56524.285
$ws.Range("H122").Value = 1924.625
$ws.Range("I122").Value = 1719.8667
$ws.Range("J122").Value = 4996
$ws.Range("K122").Value = 5159.6001
$ws.Range("L122").Value = 14988
$ws.Range("M122").Value = -2709.6001
$ws.Range("N122").Value = -19888
$ws.Range("H132").Value = 1040.8
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1040.8
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 3122.4
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -8182.4
